$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 120, shifting existing rows 120-192 down to 123-195
$ws.Rows("120:122").Insert()

# Row 120: new Early Burlat entry (Macroferia Regional de Talca, Maule)
$ws.Range("A120").Value = 5
$ws.Range("B120").Value = 'Macroferia Regional de Talca'
$ws.Range("C120").Value = 'Maule'
$ws.Range("D120").Value = 44879
$ws.Range("E120").Value = 7
$ws.Range("F120").Value = 'Fruta'
$ws.Range("G120").Value = 100103
$ws.Range("H120").Value = 'Frutos de hueso (carozo)'
$ws.Range("I120").Value = 100103001
$ws.Range("J120").Value = 'Cereza'
$ws.Range("K120").Value = 'Early Burlat'
$ws.Range("L120").Value = 'Primera'
$ws.Range("M120").Value = 50
$ws.Range("N120").Value = 25000
$ws.Range("O120").Value = 25000
$ws.Range("P120").Value = 25000
$ws.Range("Q120").Value = '$/bandeja 10 kilos'
$ws.Range("R120").Value = 'Provincia de Curicó'
$ws.Range("S120").Value = 2500
$ws.Range("T120").Value = 10

# Row 121: new Early Burlat entry (Macroferia Regional de Talca, Maule)
$ws.Range("A121").Value = 5
$ws.Range("B121").Value = 'Macroferia Regional de Talca'
$ws.Range("C121").Value = 'Maule'
$ws.Range("D121").Value = 44879
$ws.Range("E121").Value = 7
$ws.Range("F121").Value = 'Fruta'
$ws.Range("G121").Value = 100103
$ws.Range("H121").Value = 'Frutos de hueso (carozo)'
$ws.Range("I121").Value = 100103001
$ws.Range("J121").Value = 'Cereza'
$ws.Range("K121").Value = 'Early Burlat'
$ws.Range("L121").Value = 'Segunda'
$ws.Range("M121").Value = 30
$ws.Range("N121").Value = 20000
$ws.Range("O121").Value = 20000
$ws.Range("P121").Value = 20000
$ws.Range("Q121").Value = '$/bandeja 10 kilos'
$ws.Range("R121").Value = 'Provincia de Curicó'
$ws.Range("S121").Value = 2000
$ws.Range("T121").Value = 10

# Row 122: new Early Burlat entry (Macroferia Regional de Talca, Maule)
$ws.Range("A122").Value = 5
$ws.Range("B122").Value = 'Macroferia Regional de Talca'
$ws.Range("C122").Value = 'Maule'
$ws.Range("D122").Value = 44879
$ws.Range("E122").Value = 7
$ws.Range("F122").Value = 'Fruta'
$ws.Range("G122").Value = 100103
$ws.Range("H122").Value = 'Frutos de hueso (carozo)'
$ws.Range("I122").Value = 100103001
$ws.Range("J122").Value = 'Cereza'
$ws.Range("K122").Value = 'Early Burlat'
$ws.Range("L122").Value = 'Tercera'
$ws.Range("M122").Value = 15
$ws.Range("N122").Value = 15000
$ws.Range("O122").Value = 15000
$ws.Range("P122").Value = 15000
$ws.Range("Q122").Value = '$/bandeja 10 kilos'
$ws.Range("R122").Value = 'Provincia de Curicó'
$ws.Range("S122").Value = 1500
$ws.Range("T122").Value = 10
